# Apply edits described in the diff to the "Platform Coverage" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")
$ws.Activate()

# --- Cell value updates ---
$ws.Range("G2").Value = 15
$ws.Range("P2").Value = 0.8

$ws.Range("G3").Value = 50
$ws.Range("P3").Value = 0.5

$ws.Range("P4").Value = 0.5

# --- Selection / view update ---
$ws.Range("O5").Select()
$excel.ActiveWindow.ScrollColumn = $ws.Range("G1").Column
